$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column C ("Task Type"): "Non-Technical " -> "Techinical" for rows 19-30 ---
for ($r = 19; $r -le 30; $r++) {
    $ws.Cells.Item($r, 3).Value = "Techinical"
}

# --- Column H ("Status"): "Pending" -> "Done" for rows 27-30 ---
for ($r = 27; $r -le 30; $r++) {
    $ws.Cells.Item($r, 8).Value = "Done"
}

# --- New rows 31-33: three new "Update ... CDD docment According to review" tasks ---
$newRowsData = @(
    @{ Row = 31; A = "Moamen Ahmed"; B = "Update GDD docment According to review" },
    @{ Row = 32; A = "Moamen Ahmed"; B = "Update buzzer CDD docment According to review" },
    @{ Row = 33; A = "Moamen Ahmed"; B = "Update switch CDD docment According to review" }
)

foreach ($rd in $newRowsData) {
    $r = $rd.Row
    $ws.Cells.Item($r, 1).Value = $rd.A
    $ws.Cells.Item($r, 2).Value = $rd.B
    $ws.Cells.Item($r, 3).Value = "Techinical"
    $ws.Cells.Item($r, 4).Value = "3/14/2020"
    $ws.Cells.Item($r, 5).Value = "3/14/2020"
    $ws.Cells.Item($r, 6).Value = "1 day"
    $ws.Cells.Item($r, 8).Value = "Done"
}

# Row 32 and 33 have longer descriptions that wrap, so their row height grows.
$ws.Range("B32").WrapText = $true
$ws.Range("B33").WrapText = $true
$ws.Rows.Item(32).RowHeight = 30.75
$ws.Rows.Item(33).RowHeight = 31.5

# --- Column B width: 35.43 -> 37.43 (closest achievable via COM rounding) ---
$ws.Columns.Item(2).ColumnWidth = 36.59
